# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation" everywhere it appears
# - Narrow the "Status"/locale columns that previously held that long string
#   (target width ~13.41 chars; COM's ColumnWidth setter only lands on
#   1/6-character increments, so 12.5 is the input that rounds to the
#   closest achievable width, ~13.33)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E:E").ColumnWidth = 12.5
$overview.Range("F:F").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C:C").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C:C").ColumnWidth = 12.5
